$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 36.25720800000001
$ws.Range("B3").Value = 36.4254904
$ws.Range("C3").Value = 0.004641350210970341
$ws.Range("B4").Value = 41.75
$ws.Range("C4").Value = 0.1514951730425573
$ws.Range("B5").Value = 41.919
$ws.Range("C5").Value = 0.1561563151801426
$ws.Range("B6").Value = 42.11715
$ws.Range("C6").Value = 0.161621435384655
